$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "metrec" path for carol's Met_rec_comp, mirroring the existing
# entry in row 18 (column E holds matti's path already).
$ws.Range("D18").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Recycling\Met_rec_comp"

# Separate historical and projected final demand: add two new rows below
# the existing "Merged FD" row, pointing at the same Final Demand folder.
$ws.Range("A23").Value = "History"
$ws.Range("D23").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand"

$ws.Range("A24").Value = "FD"
$ws.Range("D24").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand"

$ws.Range("D23:D24").Select()
